# Update column C ("Förändrad") for data rows 2-28: 45531 -> 45532
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45531) {
        $cell.Value = 45532
    }
}
